$d = $word.ActiveDocument

function Replace-One($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $new, 1) | Out-Null
}

Replace-One "2025-07-22 Tuesday" "2025-07-23 Wednesday"
Replace-One "27÷3=9, 0" "54÷8=6, 6"
Replace-One "90÷2=45, 0" "20÷4=5, 0"
Replace-One "18÷3=6, 0" "32÷6=5, 2"
Replace-One "14÷8=1, 6" "63÷3=21, 0"
Replace-One "48÷6=8, 0" "89÷7=12, 5"
Replace-One "13÷5=2, 3" "87÷4=21, 3"
Replace-One "40÷3=13, 1" "75÷3=25, 0"
Replace-One "20÷9=2, 2" "33÷3=11, 0"
Replace-One "97÷5=19, 2" "68÷7=9, 5"
Replace-One "16÷8=2, 0" "13÷6=2, 1"
Replace-One "71÷2=35, 1" "97÷3=32, 1"
Replace-One "39÷6=6, 3" "70÷5=14, 0"
Replace-One "54÷9=6, 0" "64÷4=16, 0"
Replace-One "49÷8=6, 1" "77÷4=19, 1"
Replace-One "71÷2=35, 1" "91÷5=18, 1"
Replace-One "78÷9=8, 6" "56÷4=14, 0"
Replace-One "49÷2=24, 1" "70÷2=35, 0"
Replace-One "99÷7=14, 1" "31÷2=15, 1"
Replace-One "91÷9=10, 1" "41÷4=10, 1"
Replace-One "93÷7=13, 2" "93÷2=46, 1"
Replace-One "32÷5=6, 2" "87÷2=43, 1"
Replace-One "35÷3=11, 2" "89÷8=11, 1"
Replace-One "35÷5=7, 0" "60÷9=6, 6"
Replace-One "21÷6=3, 3" "75÷9=8, 3"
Replace-One "41÷6=6, 5" "42÷9=4, 6"
